$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Write the new/changed cell text FIRST, in the exact order the new part
#    links were authored, so the shared-string table grows in that order
#    (atmel, invensense, kemet, kingbright).
# ---------------------------------------------------------------------------

# D6 (ATMega row): new ATMEGA328P-AUR link.
$ws.Range("D6").Value = "http://www.digikey.com/product-detail/en/atmel/ATMEGA328P-AUR/ATMEGA328P-AURCT-ND/3789455"

# D3 (MPU row): was the sparkfun 13762 link, now the invensense MPU-9250 part link.
$ws.Range("D3").Value = "http://www.digikey.com/product-detail/en/invensense/MPU-9250/1428-1019-1-ND/4626450"

# C4/D4 (Capacitors row): quantity + kemet capacitor part link.
$ws.Range("C4").Value = 20
$ws.Range("D4").Value = "http://www.digikey.com/product-detail/en/kemet/C1206C106K4PACTU/399-5091-1-ND/1465625"

# D12 (red led row): kingbright LED part link.
$ws.Range("D12").Value = "http://www.digikey.com/product-detail/en/kingbright/APT3216SURCK/754-1143-1-ND/1747860"

# D2 (MPU Breakout row): the sparkfun 13762 link moves here from D3 (reuses the
# existing shared string - no new entry needed).
$ws.Range("D2").Value = "https://www.sparkfun.com/products/13762"

# ---------------------------------------------------------------------------
# 2) Rebuild the hyperlink relationships in their final left-to-right order.
#    (The engine's Hyperlinks collection is sheet-global, so the only reliable
#    way to land the unaffected links back on their original rIds - and the
#    two changed/new ones on the right rIds - is to recreate all of them here,
#    in order; each Add() reuses the text already written above.)
# ---------------------------------------------------------------------------
$ws.Range("D3").Hyperlinks.Delete()
$ws.Range("D3").Style = "Normal"

$ws.Hyperlinks.Add($ws.Range("D20"), "http://www.digikey.com/product-detail/en/vishay-semiconductor-opto-division/4N35/4N35-ND/1738522")
$ws.Range("D20").Style = "Hyperlink"
$ws.Range("D20").VerticalAlignment = -4108

$ws.Hyperlinks.Add($ws.Range("D9"), "http://www.digikey.com/product-detail/en/ftdi-future-technology-devices-international-ltd/FT232RL-REEL/768-1007-1-ND/1836402")
$ws.Range("D9").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D10"), "http://www.digikey.com/product-detail/en/kingbright/APTL3216CGCK/754-1162-1-ND/1747879")
$ws.Range("D10").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D11"), "http://www.digikey.com/product-detail/en/lite-on-inc/LTST-C150AKT/160-1166-1-ND/269238")
$ws.Range("D11").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D18"), "http://www.mouser.com/ProductDetail/AVX/12065C104KAT2A/?qs=WcAEtQi8OAaNkVk1OXh0eQ%3D%3D&gclid=CJiax9HUx9ECFcpXDQodj0kB4A")
$ws.Range("D18").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D5"), "https://www.sparkfun.com/products/8533")
$ws.Range("D5").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D2"), "https://www.sparkfun.com/products/13762")
$ws.Range("D2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D19"), "https://www.pololu.com/product/2181")
$ws.Range("D19").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D24"), "https://www.sparkfun.com/products/595")
$ws.Range("D24").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D17"), "http://www.digikey.com/product-detail/en/littelfuse-inc/0466002.NR/F1457CT-ND/521355")
$ws.Range("D17").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D6"), "http://www.digikey.com/product-detail/en/atmel/ATMEGA328P-AUR/ATMEGA328P-AURCT-ND/3789455")
$ws.Range("D6").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 3) View state: scroll/selection now rests on D12.
# ---------------------------------------------------------------------------
$ws.Range("D12").Select()
